# Update gh-pages to output generated at 456a3b4
# Updates "想去人数" (F) counts and marks one ticket tier as "不可售" (G)
# on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1168
$ws1.Range("F5").Value = 191
$ws1.Range("G5").Value = "不可售"
$ws1.Range("F6").Value = 14
$ws1.Range("F8").Value = 280
$ws1.Range("F14").Value = 164
$ws1.Range("F15").Value = 13109
$ws1.Range("F19").Value = 5388
$ws1.Range("F20").Value = 5553
$ws1.Range("F21").Value = 10

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1168
$ws4.Range("F5").Value = 191
$ws4.Range("G5").Value = "不可售"
$ws4.Range("F13").Value = 14
$ws4.Range("F24").Value = 280
$ws4.Range("F36").Value = 164
$ws4.Range("F37").Value = 13110
$ws4.Range("F42").Value = 5388
$ws4.Range("F43").Value = 5553
$ws4.Range("F44").Value = 10
